$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44764
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112012
$ws.Range("G10").Value = "Espinaca"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 35
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 15000
$ws.Range("N10").Value = "`$/cuna 10 kilos"
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 1500
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = "Hortaliza"
